$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the data that was previously in row 4, plus a new timestamp.
# Force the reference/date columns to remain text (not auto-converted to a
# number/date) the same way they were originally stored, then restore the
# default "Normal" style so no stray number-format style is left behind.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2582873016"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = "Crémone F8 Variable  L1380 1201-1600"
$ws.Range("D2").Value = "E2"
$ws.Range("E2").Value = 15

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2025-06-02 12:56:51"
$ws.Range("H2").Style = "Normal"

# Delete rows 3 and 4 entirely (remaining rows shift up)
$ws.Range("A3:H4").EntireRow.Delete()
